# Add a small "page number" text box to Slide 2, positioned in the
# bottom-right corner of the slide (bottom-right of the blue bar layout),
# showing the digit "1".
#
# Implementation notes:
#  - We clone the existing "Text 2" shape (via Copy/Paste) rather than
#    Shapes.AddShape/AddTextbox so the new shape inherits the same
#    "empty" <a:ln/> line format that the template shapes already use
#    (AddShape/AddTextbox would instead stamp a creationId extLst block
#    and require extra Fill/Line plumbing to reach the same XML shape).
#  - PowerPoint's Shape.Left/Top/Width/Height (and the AddShape args) are
#    expressed in points; the target geometry below is defined in EMU
#    (1 pt = 12700 EMU). The literal point values used are the nearest
#    values that convert back to the exact target EMU amounts.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# "Text 2" is the body-copy placeholder already on the slide; duplicate it
# as the starting point for the new page-number shape.
$src = $s.Shapes.Item(3)
$src.Copy()
$pageNum = $s.Shapes.Paste().Item(1)

$pageNum.Name = "Text 3"

# Position/size (EMU 8820000,4860000 / 216000x216000), supplied in points.
$pageNum.Left = 694.4882202148438
$pageNum.Top = 382.67718505859375
$pageNum.Width = 17.007875442504883
$pageNum.Height = 17.007875442504883

# Center the text both vertically (body anchor) and horizontally (paragraph
# alignment) inside the small square shape.
$pageNum.TextFrame.VerticalAnchor = 3

$tr = $pageNum.TextFrame.TextRange
$tr.Text = "1"
$tr.Font.Size = 8
$tr.Font.Color.RGB = 0xD9D9D9
$tr.ParagraphFormat.Alignment = 2
